$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-11-13 -> 2023-11-14) for every data row (rows 2-56).
$ws.Range("C2:C56").Value = 45244
